$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Andorra): language D -> E, nationality andorranisch -> Andorran ---
$ws.Range("B2").Value = "E"
$ws.Range("E2").Value = "Andorran"
$ws.Range("G2").Value = "Andorran"

# --- Row 3: was AE/UAE data, now becomes BE/Belgium data ---
$ws.Range("B3").Value = "E"
$ws.Range("C3").Value = "BE"
$ws.Range("D3").Value = "Belgium"
$ws.Range("E3").Value = "Belgian"
$ws.Range("F3").Value = "Belgium"
$ws.Range("G3").Value = "Belgian"

# --- Row 4 (new): Germany ---
# Copy A2 ("001" stored as text) down so the leading-zero code isn't
# reinterpreted as the number 1.
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B4").Value = "E"
$ws.Range("C4").Value = "DE"
$ws.Range("D4").Value = "Germany"
$ws.Range("E4").Value = "German"
$ws.Range("F4").Value = "Germany"
$ws.Range("G4").Value = "German"

# --- Row 5 (new): Micronesia ---
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("B5").Value = "E"
$ws.Range("C5").Value = "FM"
$ws.Range("D5").Value = "Micronesia"
$ws.Range("E5").Value = "Micronesian"
$ws.Range("F5").Value = "Micronesia"
$ws.Range("G5").Value = "Micronesian"

# --- Re-apply the AutoFilter over the grown range, keeping only "AD" visible ---
# Clear the old (now too-small) AutoFilter range first so re-applying actually
# grows <autoFilter ref="..."> to cover the new rows instead of leaving it
# clamped to the original A1:G3.
$ws.AutoFilterMode = $false
# Operator 7 = xlFilterValues, feeding an explicit list of values to keep.
$ws.Range("A1:G5").AutoFilter(3, @("AD"), 7)

# --- Keep the hidden _FilterDatabase defined name in sync with the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Internal table!_FilterDatabase") {
        $n.RefersTo = "='Internal table'!`$A`$1:`$G`$5"
    }
}
